# "maj template comment a la fin"
#
# Move the "Comment" column to the end of the table (after
# AmpliconSizeUnit), shifting the primer/amplicon columns that used to
# follow it one position to the left, for every row of the template
# (the 4 header/description rows plus the example row).
#
# Before:  ... SampleID | Comment | ForwardPrimerName | ... | AmpliconSizeUnit
# After:   ... SampleID | ForwardPrimerName | ... | AmpliconSizeUnit | Comment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow  = $firstRow + $used.Rows.Count - 1
$lastCol  = $firstCol + $used.Columns.Count - 1

# Column "Comment" is the 10th used column (A..I = 9 columns before it),
# i.e. column J when the table starts at column A.
$commentCol = $firstCol + 9

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Capture the row's current values from the Comment column through the
    # last column (Comment, ForwardPrimerName, ..., AmpliconSizeUnit).
    $values = @()
    for ($c = $commentCol; $c -le $lastCol; $c++) {
        $values += , ($ws.Cells.Item($r, $c).Value2)
    }

    # Shift everything after Comment one column to the left.
    for ($c = $commentCol; $c -le ($lastCol - 1); $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - $commentCol + 1]
    }

    # Place the original Comment value in the now-vacant last column.
    $ws.Cells.Item($r, $lastCol).Value = $values[0]
}
